$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of an existing header cell (G1) onto the new header cell H1,
# so it reuses the same header style (bold, centered, bordered) instead of
# creating a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
